# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (interest count) values in column F for a handful
# of rows on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8
$ws1.Range("F5").Value = 79
$ws1.Range("F6").Value = 5248
$ws1.Range("F7").Value = 177
$ws1.Range("F8").Value = 87
$ws1.Range("F9").Value = 100
$ws1.Range("F10").Value = 354
$ws1.Range("F11").Value = 13

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8
$ws4.Range("F9").Value = 79
$ws4.Range("F10").Value = 5248
$ws4.Range("F11").Value = 177
$ws4.Range("F12").Value = 87
$ws4.Range("F13").Value = 100
$ws4.Range("F15").Value = 354
$ws4.Range("F16").Value = 13
